# Generate Report for Handoff
# Update status + timestamps on each localization sheet, then resize the
# "Status"-related columns to fit their (now shorter) new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E, F) + HO generate date (G) ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 02:59:13"

# --- zh-cn sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 02:59:09"

# --- de-de sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 02:59:13"

# Resize the status columns to fit the new, shorter text (autofit-style
# shrink: "Ready for handoff" is narrower than the old status string).
# 16.34 is the input that the engine's character->pixel quantization maps
# to the closest reproducible width (matches the target ~17.22 char width).
$wsOverview.Columns.Item(5).ColumnWidth = 16.34
$wsOverview.Columns.Item(6).ColumnWidth = 16.34
$wsZhCn.Columns.Item(3).ColumnWidth = 16.34
$wsDeDe.Columns.Item(3).ColumnWidth = 16.34
